# Update NATMI LR-pair TPM-derived values (Efna2-Epha7) for sheet1
# as recomputed after switching to new TPM data (commit: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 3.034777666666667
$ws.Range("H2").Value2 = 9.104333
$ws.Range("I2").Value2 = 0.2502264227183869
$ws.Range("J2").Value2 = 0.2502264227183869
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 0.042868
$ws.Range("N2").Value2 = 0.128604
$ws.Range("O2").Value2 = 0.03014606792405771
$ws.Range("P2").Value2 = 0.03014606792405771
$ws.Range("Q2").Value2 = 0.1300948490146667
$ws.Range("R2").Value2 = 1.170853641132
$ws.Range("S2").Value2 = 0.007543342735662469
$ws.Range("T2").Value2 = 0.00754334273566247
# Row 3
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 3.034777666666667
$ws.Range("H3").Value2 = 9.104333
$ws.Range("I3").Value2 = 0.2502264227183869
$ws.Range("J3").Value2 = 0.2502264227183869
$ws.Range("O3").Value2 = 0.2718481285523376
$ws.Range("P3").Value2 = 0.2718481285523376
$ws.Range("Q3").Value2 = 1.173156025788445
$ws.Range("R3").Value2 = 10.558404232096
$ws.Range("S3").Value2 = 0.06802358473033961
$ws.Range("T3").Value2 = 0.06802358473033961
# Row 4
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 3.034777666666667
$ws.Range("H4").Value2 = 9.104333
$ws.Range("I4").Value2 = 0.2502264227183869
$ws.Range("J4").Value2 = 0.2502264227183869
$ws.Range("M4").Value2 = 0.9839956666666666
$ws.Range("N4").Value2 = 2.951987
$ws.Range("O4").Value2 = 0.69197537100662
$ws.Range("P4").Value2 = 0.69197537100662
$ws.Range("Q4").Value2 = 2.986208073296778
$ws.Range("R4").Value2 = 26.875872659671
$ws.Range("S4").Value2 = 0.1731505216962151
$ws.Range("T4").Value2 = 0.1731505216962151
# Row 5
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 3.034777666666667
$ws.Range("H5").Value2 = 9.104333
$ws.Range("I5").Value2 = 0.2502264227183869
$ws.Range("J5").Value2 = 0.2502264227183869
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.008575333333333332
$ws.Range("N5").Value2 = 0.025726
$ws.Range("O5").Value2 = 0.006030432516984765
$ws.Range("P5").Value2 = 0.006030432516984765
$ws.Range("Q5").Value2 = 0.02602423008422222
$ws.Range("R5").Value2 = 0.234218070758
$ws.Range("S5").Value2 = 0.001508973556169736
$ws.Range("T5").Value2 = 0.001508973556169736
# Row 6
$ws.Range("I6").Value2 = 0.4835045831069426
$ws.Range("J6").Value2 = 0.4835045831069426
$ws.Range("K6").Value2 = 1
$ws.Range("L6").Value2 = 0.3333333333333333
$ws.Range("M6").Value2 = 0.042868
$ws.Range("N6").Value2 = 0.128604
$ws.Range("O6").Value2 = 0.03014606792405771
$ws.Range("P6").Value2 = 0.03014606792405771
$ws.Range("Q6").Value2 = 0.2513781520506666
$ws.Range("R6").Value2 = 2.262403368456
$ws.Range("S6").Value2 = 0.0145757620039351
$ws.Range("T6").Value2 = 0.0145757620039351
# Row 7
$ws.Range("I7").Value2 = 0.4835045831069426
$ws.Range("J7").Value2 = 0.4835045831069426
$ws.Range("O7").Value2 = 0.2718481285523376
$ws.Range("P7").Value2 = 0.2718481285523376
$ws.Range("S7").Value2 = 0.1314398160641005
$ws.Range("T7").Value2 = 0.1314398160641005
# Row 8
$ws.Range("I8").Value2 = 0.4835045831069426
$ws.Range("J8").Value2 = 0.4835045831069426
$ws.Range("M8").Value2 = 0.9839956666666666
$ws.Range("N8").Value2 = 2.951987
$ws.Range("O8").Value2 = 0.69197537100662
$ws.Range("P8").Value2 = 0.69197537100662
$ws.Range("Q8").Value2 = 5.77015518131311
$ws.Range("R8").Value2 = 51.931396631818
$ws.Range("S8").Value2 = 0.3345732632788277
$ws.Range("T8").Value2 = 0.3345732632788277
# Row 9
$ws.Range("I9").Value2 = 0.4835045831069426
$ws.Range("J9").Value2 = 0.4835045831069426
$ws.Range("K9").Value2 = 1
$ws.Range("L9").Value2 = 0.3333333333333333
$ws.Range("M9").Value2 = 0.008575333333333332
$ws.Range("N9").Value2 = 0.025726
$ws.Range("O9").Value2 = 0.006030432516984765
$ws.Range("P9").Value2 = 0.006030432516984765
$ws.Range("Q9").Value2 = 0.05028579468488888
$ws.Range("R9").Value2 = 0.4525721521639999
$ws.Range("S9").Value2 = 0.002915741760079269
$ws.Range("T9").Value2 = 0.002915741760079269
# Row 10
$ws.Range("G10").Value2 = 2.564975
$ws.Range("H10").Value2 = 7.694925
$ws.Range("I10").Value2 = 0.2114897989601526
$ws.Range("J10").Value2 = 0.2114897989601526
$ws.Range("K10").Value2 = 1
$ws.Range("L10").Value2 = 0.3333333333333333
$ws.Range("M10").Value2 = 0.042868
$ws.Range("N10").Value2 = 0.128604
$ws.Range("O10").Value2 = 0.03014606792405771
$ws.Range("P10").Value2 = 0.03014606792405771
$ws.Range("Q10").Value2 = 0.1099553483
$ws.Range("R10").Value2 = 0.9895981346999999
$ws.Range("S10").Value2 = 0.006375585844698071
$ws.Range("T10").Value2 = 0.006375585844698071
# Row 11
$ws.Range("G11").Value2 = 2.564975
$ws.Range("H11").Value2 = 7.694925
$ws.Range("I11").Value2 = 0.2114897989601526
$ws.Range("J11").Value2 = 0.2114897989601526
$ws.Range("O11").Value2 = 0.2718481285523376
$ws.Range("P11").Value2 = 0.2718481285523376
$ws.Range("Q11").Value2 = 0.9915440957333334
$ws.Range("R11").Value2 = 8.923896861599999
$ws.Range("S11").Value2 = 0.0574931060552276
$ws.Range("T11").Value2 = 0.0574931060552276
# Row 12
$ws.Range("G12").Value2 = 2.564975
$ws.Range("H12").Value2 = 7.694925
$ws.Range("I12").Value2 = 0.2114897989601526
$ws.Range("J12").Value2 = 0.2114897989601526
$ws.Range("M12").Value2 = 0.9839956666666666
$ws.Range("N12").Value2 = 2.951987
$ws.Range("O12").Value2 = 0.69197537100662
$ws.Range("P12").Value2 = 0.69197537100662
$ws.Range("Q12").Value2 = 2.523924285108333
$ws.Range("R12").Value2 = 22.715318565975
$ws.Range("S12").Value2 = 0.1463457320995671
$ws.Range("T12").Value2 = 0.1463457320995671
# Row 13
$ws.Range("G13").Value2 = 2.564975
$ws.Range("H13").Value2 = 7.694925
$ws.Range("I13").Value2 = 0.2114897989601526
$ws.Range("J13").Value2 = 0.2114897989601526
$ws.Range("K13").Value2 = 1
$ws.Range("L13").Value2 = 0.3333333333333333
$ws.Range("M13").Value2 = 0.008575333333333332
$ws.Range("N13").Value2 = 0.025726
$ws.Range("O13").Value2 = 0.006030432516984765
$ws.Range("P13").Value2 = 0.006030432516984765
$ws.Range("Q13").Value2 = 0.02199551561666667
$ws.Range("R13").Value2 = 0.19795964055
$ws.Range("S13").Value2 = 0.001275374960659875
$ws.Range("T13").Value2 = 0.001275374960659875
# Row 14
$ws.Range("E14").Value2 = 3
$ws.Range("F14").Value2 = 1
$ws.Range("G14").Value2 = 0.6643690000000001
$ws.Range("H14").Value2 = 1.993107
$ws.Range("I14").Value2 = 0.05477919521451775
$ws.Range("J14").Value2 = 0.05477919521451775
$ws.Range("K14").Value2 = 1
$ws.Range("L14").Value2 = 0.3333333333333333
$ws.Range("M14").Value2 = 0.042868
$ws.Range("N14").Value2 = 0.128604
$ws.Range("O14").Value2 = 0.03014606792405771
$ws.Range("P14").Value2 = 0.03014606792405771
$ws.Range("Q14").Value2 = 0.028480170292
$ws.Range("R14").Value2 = 0.256321532628
$ws.Range("S14").Value2 = 0.001651377339762069
$ws.Range("T14").Value2 = 0.001651377339762069
# Row 15
$ws.Range("E15").Value2 = 3
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 0.6643690000000001
$ws.Range("H15").Value2 = 1.993107
$ws.Range("I15").Value2 = 0.05477919521451775
$ws.Range("J15").Value2 = 0.05477919521451775
$ws.Range("O15").Value2 = 0.2718481285523376
$ws.Range("P15").Value2 = 0.2718481285523376
$ws.Range("Q15").Value2 = 0.2568255672426667
$ws.Range("R15").Value2 = 2.311430105184
$ws.Range("S15").Value2 = 0.01489162170266982
$ws.Range("T15").Value2 = 0.01489162170266982
# Row 16
$ws.Range("E16").Value2 = 3
$ws.Range("F16").Value2 = 1
$ws.Range("G16").Value2 = 0.6643690000000001
$ws.Range("H16").Value2 = 1.993107
$ws.Range("I16").Value2 = 0.05477919521451775
$ws.Range("J16").Value2 = 0.05477919521451775
$ws.Range("M16").Value2 = 0.9839956666666666
$ws.Range("N16").Value2 = 2.951987
$ws.Range("O16").Value2 = 0.69197537100662
$ws.Range("P16").Value2 = 0.69197537100662
$ws.Range("Q16").Value2 = 0.6537362170676667
$ws.Range("R16").Value2 = 5.883625953609
$ws.Range("S16").Value2 = 0.03790585393200999
$ws.Range("T16").Value2 = 0.03790585393200998
# Row 17
$ws.Range("E17").Value2 = 3
$ws.Range("F17").Value2 = 1
$ws.Range("G17").Value2 = 0.6643690000000001
$ws.Range("H17").Value2 = 1.993107
$ws.Range("I17").Value2 = 0.05477919521451775
$ws.Range("J17").Value2 = 0.05477919521451775
$ws.Range("K17").Value2 = 1
$ws.Range("L17").Value2 = 0.3333333333333333
$ws.Range("M17").Value2 = 0.008575333333333332
$ws.Range("N17").Value2 = 0.025726
$ws.Range("O17").Value2 = 0.006030432516984765
$ws.Range("P17").Value2 = 0.006030432516984765
$ws.Range("Q17").Value2 = 0.005697185631333333
$ws.Range("R17").Value2 = 0.051274670682
$ws.Range("S17").Value2 = 0.000330342240075884
$ws.Range("T17").Value2 = 0.000330342240075884
